$wb = $excel.ActiveWorkbook

# Sheet1: RO Non-Availability Dates - insert new column F for White River Junction (RO05)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Columns.Item(6).Insert()
$ws1.Columns.Item(6).ColumnWidth = 14.666666666666666
$ws1.Range("F3").Value = "White River Junction, VT"
$ws1.Range("F2").Value = "RO05"

# Sheet3: RO Allocations - insert new row 7 for White River Junction (RO05)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Rows.Item(7).Insert()
$ws3.Range("A8:N8").Copy()
$ws3.Range("A7:N7").PasteSpecial(-4122)
$ws3.Range("B7").Value = "White River Junction, VT"
$ws3.Range("C7").Value = "RO05"
$ws3.Range("D7").Value = 0
$ws3.Range("E7").Value = 0
$ws3.Range("F7").Value = 0
$ws3.Range("G7").Value = 0
$ws3.Range("H7").Value = 0
